# cane_composition_data.xlsx edit:
#  - split off the last data row ("EC") of the Summarized sheet into a new,
#    unused sheet ("NotUsed") while fixing the yield label typo
#    (MT/hc -> MT/ha) on the "Raw data" sheet.

$wb = $excel.ActiveWorkbook

$wsSummarized = $wb.Worksheets.Item("Summarized")
$wsRawData    = $wb.Worksheets.Item("Raw data")

# ---------------------------------------------------------------------
# 1) Fix the typo in the yield label used on the "Raw data" sheet
#    ("Yield (dry MT/hc)" -> "Yield (dry MT/ha)").
# ---------------------------------------------------------------------
$wsRawData.Range("F24").Value = "Yield (dry MT/ha)"

# ---------------------------------------------------------------------
# 2) Re-merge G12:H12 on "Raw data" (moves it to the end of the
#    worksheet's merged-cell list).
# ---------------------------------------------------------------------
$wsRawData.Range("G12:H12").UnMerge()
$wsRawData.Range("G12:H12").Merge()

# ---------------------------------------------------------------------
# 3) Create the new "NotUsed" sheet at the end of the workbook and move
#    the "EC" row (row 15) of "Summarized" into it as row 1, preserving
#    formulas/values/number formats.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNotUsed = $wb.Worksheets.Add($null, $lastSheet)
$wsNotUsed.Name = "NotUsed"

$wsNotUsed.Range("A1").Value = "EC"

$wsNotUsed.Range("B1").Formula = "=1-D1-F1-0.07"
$wsNotUsed.Range("B1").NumberFormat = "0.00%"

$wsNotUsed.Range("C1").Value = 0.001
$wsNotUsed.Range("C1").NumberFormat = "0.00%"

$wsNotUsed.Range("D1").Formula = "=Summarized!D6"
$wsNotUsed.Range("D1").NumberFormat = "0.00%"

$wsNotUsed.Range("E1").Value = 0.000001
$wsNotUsed.Range("E1").NumberFormat = "0.00%"

$wsNotUsed.Range("F1").Formula = "=0.091/0.4"
$wsNotUsed.Range("F1").NumberFormat = "0.00%"

$wsNotUsed.Range("G1").Value = 0.001
$wsNotUsed.Range("G1").NumberFormat = "0.00%"

$wsNotUsed.Range("H1").Formula = "=0.6"
$wsNotUsed.Range("H1").NumberFormat = "0.00%"

$wsNotUsed.Range("I1").Formula = "='Raw data'!D23/100"
$wsNotUsed.Range("I1").NumberFormat = "0.00%"

$wsNotUsed.Range("J1").Formula = "=L1*'Raw data'!`$F`$25"
$wsNotUsed.Range("J1").NumberFormat = "0.00"

$wsNotUsed.Range("K1").Formula = "=M1*'Raw data'!`$F`$25"
$wsNotUsed.Range("K1").NumberFormat = "0.00"

$wsNotUsed.Range("L1").Formula = "='Raw data'!F27/'Raw data'!F25"
$wsNotUsed.Range("L1").NumberFormat = "0.00%"

$wsNotUsed.Range("M1").Formula = "='Raw data'!L23/100"
$wsNotUsed.Range("M1").NumberFormat = "0.00%"

$wsNotUsed.Range("A1:M1").Select()

# Now that the row lives on "NotUsed", delete it from "Summarized".
$wsSummarized.Rows.Item(15).Delete()

# ---------------------------------------------------------------------
# 4) Restore the view state: "Summarized" stays the active sheet/tab,
#    with the same row selected that was left selected after the row
#    was cut out.
# ---------------------------------------------------------------------
$wsSummarized.Activate()
$wsSummarized.Range("A13:M13").Select()
